$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.208814144134521
$ws.Range("B1").Value = 4.430998802185059
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 4.123380661010742
$ws.Range("E1").Value = 2.166831493377686
